# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 01878cfb-... handback row (row 18) on both the
# "zh-cn" and "de-de" report sheets, regenerating the handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D18").Value = "2016-03-07 02:11:06"
$wsZhCn.Range("G18").Value = "2016-03-07 02:11:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D18").Value = "2016-03-07 02:11:15"
$wsDeDe.Range("G18").Value = "2016-03-07 02:12:10"
